$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the ergm description table contents to match the revised terminology.
$ws.Range('A1').Value2 = 'Network Property'
$ws.Range('B1').Value2 = 'Configurations (effect/parameter)'
$ws.Range('C1').Value2 = 'Description'
$ws.Range('D1').Value2 = 'Interpretation'

$ws.Range('A2').Value2 = 'Density'
$ws.Range('B2').Value2 = 'Edges'
$ws.Range('C2').Value2 = 'Number of ties in the network'
$ws.Range('D2').Value2 = 'two burials having relationship'

$ws.Range('A3').Value2 = 'Node covariate of age'
$ws.Range('B3').Value2 = 'Homophily/nodematch.age '
$ws.Range('C3').Value2 = 'Density of ties between nodes with the same age'
$ws.Range('D3').Value2 = 'burials having the same age to be connected'

$ws.Range('A4').Value2 = 'Node covariate of sex'
$ws.Range('B4').Value2 = 'Homophily/nodematch.sex '
$ws.Range('C4').Value2 = 'Density of ties between nodes with the same gender'
$ws.Range('D4').Value2 = 'burials having the same gender to be connected'

$ws.Range('A5').Value2 = 'Node covariate of ritual pottery'
$ws.Range('B5').Value2 = 'Homophily/nodematch.ritual'
$ws.Range('C5').Value2 = 'Density of ties between nodes with same ritual treatment'
$ws.Range('D5').Value2 = 'burials having the same ritual treatment to be connected'

$ws.Range('A6').Value2 = 'Node covariate of burial value'
$ws.Range('B6').Value2 = 'Homophily/nodematch.value '
$ws.Range('C6').Value2 = 'Density of ties between nodes with same scale of wealth'
$ws.Range('D6').Value2 = 'burials having the same scale of wealth to be connected'

$ws.Range('A7').Value2 = 'Transitivity or cohesion'
$ws.Range('B7').Value2 = 'Geometrically weighted edgewise shared partner (gwesp)'
$ws.Range('C7').Value2 = 'Tendency for nodes with shared partners to be tied '
$ws.Range('D7').Value2 = 'burials to be connected with a third shared burial '

$ws.Range('A8').Value2 = 'Popularity'
$ws.Range('B8').Value2 = 'geometrically weighted degree distribution (gwdeg)'
$ws.Range('C8').Value2 = 'Tendency towards centralization in distribution'
$ws.Range('D8').Value2 = 'burials being connected with multiple partners'

$ws.Range('A9').Value2 = 'physical distance'
$ws.Range('B9').Value2 = 'dyadic relationship/dyadcov.distance'
$ws.Range('C9').Value2 = 'Distance (in meter) between each pairs of burials'
$ws.Range('D9').Value2 = 'burials with shorter distance to be connected based on kinship-based relations '

# Row heights follow from the wrapped-text autofit for the new content.
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 34
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 34
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 51
$ws.Rows.Item(7).RowHeight = 34
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 51

# Restore the cursor/selection to the cell shown in the saved file.
$ws.Range('A9').Select()
